$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows below the existing "R1-R6" resistor row (row 9),
# inheriting that row's formatting/styles.
$ws.Rows("10:11").Insert() | Out-Null

# Row 9: R1-R3 (was R1-R6, now split into three separate resistor lines)
$ws.Range("A9").Value = "R1-R3"
$ws.Range("B9").Value = "0603"
$ws.Range("C9").Value = "4.7k Ohm"
$ws.Range("D9").Value = "4.7k Ohm"
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = "CRCW06034K70JNEA"

# Row 10: R4-R5 (re-uses the old 1k Ohm value/part)
$ws.Range("A10").Value = "R4-R5"
$ws.Range("B10").Value = "0603"
$ws.Range("C10").Value = "1k Ohm"
$ws.Range("D10").Value = "1k Ohm "
$ws.Range("E10").Value = 2
$ws.Range("G10").Value = "CRCW06031K00FKEA"
$ws.Hyperlinks.Add($ws.Range("H10"), "http://www.vishay.com/docs/20035/dcrcwe3.pdf", "", "", "data") | Out-Null

# Row 11: R6 (new 10k Ohm resistor)
$ws.Range("A11").Value = "R6"
$ws.Range("B11").Value = "0603"
$ws.Range("C11").Value = "10k Ohm"
$ws.Range("D11").Value = "10k Ohm"
$ws.Range("E11").Value = 1
$ws.Range("G11").Value = "CRCW060310K0FKEA"
$ws.Hyperlinks.Add($ws.Range("H11"), "http://www.vishay.com/docs/20035/dcrcwe3.pdf", "", "", "data") | Out-Null

# Update the saved selection to match the author's final cursor position
$ws.Range("D10").Select() | Out-Null
